$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 64; $row++) {
    $ws.Cells.Item($row, 15).Value = "2022-09-08 21:00:56"
}
